$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 8's data)
$ws.Range("D2").Value = 44187
$ws.Range("K2").Value = "Dina"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 55
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15455
$ws.Range("Q2").Value = "$/caja 15 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1030
$ws.Range("T2").Value = 15

# Row 3 (was row 5's data)
$ws.Range("D3").Value = 44174
$ws.Range("K3").Value = "Castle Brite"
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 75
$ws.Range("N3").Value = 9000
$ws.Range("O3").Value = 10000
$ws.Range("P3").Value = 9467
$ws.Range("Q3").Value = "$/caja 10 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 947
$ws.Range("T3").Value = 10

# Row 4 (was row 7's data)
$ws.Range("D4").Value = 44165
$ws.Range("K4").Value = "Castle Brite"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 16000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 16500
$ws.Range("Q4").Value = "$/caja 15 kilos granel"
$ws.Range("R4").Value = "Provincia de Limarí"
$ws.Range("S4").Value = 1100
$ws.Range("T4").Value = 15

# Row 5 (was row 4's data)
$ws.Range("D5").Value = 44176
$ws.Range("K5").Value = "Castle Brite"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 17400
$ws.Range("Q5").Value = "$/caja 18 kilos"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 967
$ws.Range("T5").Value = 18

# Row 6 (was row 2's data)
$ws.Range("D6").Value = 44189
$ws.Range("K6").Value = "Dina"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 80
$ws.Range("N6").Value = 16000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 16562
$ws.Range("Q6").Value = "$/caja 18 kilos"
$ws.Range("R6").Value = "Región de O'Higgins"
$ws.Range("S6").Value = 920
$ws.Range("T6").Value = 18

# Row 7 (was row 6's data)
$ws.Range("D7").Value = 44168
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 17000
$ws.Range("P7").Value = 16500
$ws.Range("Q7").Value = "$/caja 16 kilos granel"
$ws.Range("R7").Value = "Región de Coquimbo"
$ws.Range("S7").Value = 1031
$ws.Range("T7").Value = 16

# Row 8 (was row 3's data)
$ws.Range("D8").Value = 44181
$ws.Range("K8").Value = "Modesto"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20500
$ws.Range("Q8").Value = "$/caja 18 kilos"
$ws.Range("R8").Value = "Región de Coquimbo"
$ws.Range("S8").Value = 1139
$ws.Range("T8").Value = 18
